$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Résultats")

# Clear the "Bénéfice optimal" value (B1) so the sheet is no longer filled in
$ws.Range("B1").ClearContents()

# Clear the "Nombre chargé" row values (B3:U3) for the same reason
$ws.Range("B3:U3").ClearContents()
